# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45990

$ws.Range("B2").Value = 105.72
$ws.Range("C2").Value = 99.53
$ws.Range("D2").Value = 97.06
$ws.Range("E2").Value = 94.38
$ws.Range("F2").Value = 84.64
$ws.Range("G2").Value = 83.26000000000001
$ws.Range("H2").Value = 85.59999999999999
$ws.Range("I2").Value = 91.92
$ws.Range("J2").Value = 93.8
$ws.Range("K2").Value = 80.56
$ws.Range("L2").Value = 62.88
$ws.Range("M2").Value = 50.06
$ws.Range("N2").Value = 50.53
$ws.Range("O2").Value = 44.4
$ws.Range("P2").Value = 47.64
$ws.Range("Q2").Value = 56.85
$ws.Range("R2").Value = 72.33
$ws.Range("S2").Value = 95.84999999999999
$ws.Range("T2").Value = 97.78
$ws.Range("U2").Value = 100.24
$ws.Range("V2").Value = 100.56
$ws.Range("W2").Value = 99.61
$ws.Range("X2").Value = 91.90000000000001
$ws.Range("Y2").Value = 86.59
$ws.Range("Z2").Value = 82.23999999999999

$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 99.17
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 102.62
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 100.08
$ws.Range("AG2").Value = "9h-16h"
